# Insert a new data row into the "Pimiento" price sheet (Vega Monumental
# Concepcion, Bio-Bio) right before the existing row that is currently on
# worksheet row 223. Every row from the old row 223 down to the old last
# row (316) shifts down by one, and the new row is populated with a fresh
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 223:316 down to 224:317 by inserting a blank row at 223.
$ws.Rows(223).Insert()

# Populate the newly inserted row 223 with the new observation.
$ws.Cells.Item(223, 1).Value = 11
$ws.Cells.Item(223, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(223, 3).Value = "Bíobío"
$ws.Cells.Item(223, 4).Value = 44726
$ws.Cells.Item(223, 5).Value = 8
$ws.Cells.Item(223, 6).Value = 100112002
$ws.Cells.Item(223, 7).Value = "Pimiento"
$ws.Cells.Item(223, 8).Value = "Morrón rojo"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 140
$ws.Cells.Item(223, 11).Value = 22000
$ws.Cells.Item(223, 12).Value = 23000
$ws.Cells.Item(223, 13).Value = 22571
$ws.Cells.Item(223, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(223, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(223, 16).Value = 1254
$ws.Cells.Item(223, 17).Value = 18
$ws.Cells.Item(223, 18).Value = "Hortaliza"
